$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "assignee"
$ws.Range("B1").Value = "reward"
$ws.Range("C1").Value = "permit1"
$ws.Range("D1").Value = "permit2"

$ws.Range("A2").Value = "devpanther"
$ws.Range("B2").Value = 159.15
$ws.Range("C2").Value = "https://pay.ubq.fi?claim=eyJwZXJtaXQiOnsicGVybWl0dGVkIjp7InRva2VuIjoiMHhlOTFEMTUzRTBiNDE1MThBMkNlOERkM0Q3OTQ0RmE4NjM0NjNhOTdkIiwiYW1vdW50IjoiMTUwMDAwMDAwMDAwMDAwMDAwMDAwIn0sIm5vbmNlIjoiNTk5NDM2NjEwMjQxMzM0NDI5MjYzNjE2ODYxMTA5MDc4NjEwNTk2MzQxOTQxMzYzMzQ0MzM0OTU0NTY1NjYxMjQwNjMwNTc4NDg1OTUiLCJkZWFkbGluZSI6IjExNTc5MjA4OTIzNzMxNjE5NTQyMzU3MDk4NTAwODY4NzkwNzg1MzI2OTk4NDY2NTY0MDU2NDAzOTQ1NzU4NDAwNzkxMzEyOTYzOTkzNSJ9LCJ0cmFuc2ZlckRldGFpbHMiOnsidG8iOiIweGY3NkYxQUNCNjYwMjBmODkzYzk1MzcxZjc0MDU0OUYzMTJERUEzZjEiLCJyZXF1ZXN0ZWRBbW91bnQiOiIxNTAwMDAwMDAwMDAwMDAwMDAwMDAifSwib3duZXIiOiIweGY4N2NhNDU4M0M3OTIyMTJlNTI3MjBkMTI3RTdFMEEzOEI4MThhRDEiLCJzaWduYXR1cmUiOiIweDBiMjQyNDg2N2FmOGU0YTM0NGExZjU4YWMwZTkxMjY3MDFhODliZjBmY2RiOTVlMTlmNzUxZmViOWVlODVmNzQ2MzY0ZjE0NzMzNTdmMGQ3ZDExNjkxZjVlN2NjNzM4MTM3NzJjYTJkNzNhZTcxZjk2NWRlZTY2YTgwZGJhZDI2MWIifQ==&network=100"
$ws.Range("D2").Value = "https://pay.ubq.fi?claim=eyJwZXJtaXQiOnsicGVybWl0dGVkIjp7InRva2VuIjoiMHhlOTFEMTUzRTBiNDE1MThBMkNlOERkM0Q3OTQ0RmE4NjM0NjNhOTdkIiwiYW1vdW50IjoiOTE1MDAwMDAwMDAwMDAwMDAwMCJ9LCJub25jZSI6IjQ1MDg0OTkyMjcxODgwMjM1MDg5NDg0MTE3ODE3MDk4MDk2MzU5MDIzOTI0MTQ1NTA2MDM0OTk3MzM3MDQ0NzIyNTExMzM2NTU2NjY4IiwiZGVhZGxpbmUiOiIxMTU3OTIwODkyMzczMTYxOTU0MjM1NzA5ODUwMDg2ODc5MDc4NTMyNjk5ODQ2NjU2NDA1NjQwMzk0NTc1ODQwMDc5MTMxMjk2Mzk5MzUifSwidHJhbnNmZXJEZXRhaWxzIjp7InRvIjoiMHhmNzZGMUFDQjY2MDIwZjg5M2M5NTM3MWY3NDA1NDlGMzEyREVBM2YxIiwicmVxdWVzdGVkQW1vdW50IjoiOTE1MDAwMDAwMDAwMDAwMDAwMCJ9LCJvd25lciI6IjB4Zjg3Y2E0NTgzQzc5MjIxMmU1MjcyMGQxMjdFN0UwQTM4QjgxOGFEMSIsInNpZ25hdHVyZSI6IjB4YzRjOWEzMWM4OTFjODZmNzRhOTExNzBhZDBkYzQ3YjNjMGU3OTQxNTIwYTQ0MDBjMTg4Zjc4YmJjMzM3NzIwZjAxZDM3OTE1ZDZjNDFkOTdkNzFhYjVmYjZmZjNlZTU4MTlmYjk1YWNjYzA5OWMyNTViNTBiNDBlM2FkNjJiNDQxYyJ9&network=100"

$ws.Range("A3").Value = "pavlovcik"
$ws.Range("B3").Value = 27.3
$ws.Range("C3").Value = "https://pay.ubq.fi?claim=eyJwZXJtaXQiOnsicGVybWl0dGVkIjp7InRva2VuIjoiMHhlOTFEMTUzRTBiNDE1MThBMkNlOERkM0Q3OTQ0RmE4NjM0NjNhOTdkIiwiYW1vdW50IjoiMjczMDAwMDAwMDAwMDAwMDAwMDAifSwibm9uY2UiOiIyNjAwNzgzNzczNTE2ODg2MTUxMTA3ODM3MjA1OTcyNjc4OTEzODgxMDExMDM4ODc3NzU0NzAwOTUxMzYxMTkwNjg1ODIyODYwMDU1MiIsImRlYWRsaW5lIjoiMTE1NzkyMDg5MjM3MzE2MTk1NDIzNTcwOTg1MDA4Njg3OTA3ODUzMjY5OTg0NjY1NjQwNTY0MDM5NDU3NTg0MDA3OTEzMTI5NjM5OTM1In0sInRyYW5zZmVyRGV0YWlscyI6eyJ0byI6IjB4NDAwN0NFMjA4M2M3RjNFMTgwOTdhZUIzQTM5YmI4ZUMxNDlhMzQxZCIsInJlcXVlc3RlZEFtb3VudCI6IjI3MzAwMDAwMDAwMDAwMDAwMDAwIn0sIm93bmVyIjoiMHhmODdjYTQ1ODNDNzkyMjEyZTUyNzIwZDEyN0U3RTBBMzhCODE4YUQxIiwic2lnbmF0dXJlIjoiMHhjNmJmMGQyYmY1YjU1MTYwODFkYjE4NzgxNmFlODFiN2E3MjgzYmQ5YzdhMjhiNDZlZDUzNzQxN2FlMmUzMGEzMTUzODNhYWFmNTQ0NmZiODYxOGVlZjQxMmRhNGFmM2ViNDAyYzMwNTU1MGRhMTNlOGY4YmE5MGFkZjc0MTE1OTFjIn0=&network=100"
